# Insert two new weekly data rows (Ají, "Americana (o)" and "Inferno") right
# before the existing row 275, shifting the remaining rows (old 275..296)
# down to 277..298, and fill the two new rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 275 (pushes everything at/after 275 down by 2)
$ws.Rows.Item(275).Insert()
$ws.Rows.Item(275).Insert()

# --- New row 275 ---
$ws.Cells.Item(275, 1).Value = 9
$ws.Cells.Item(275, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(275, 3).Value = "Metropolitana"
$ws.Cells.Item(275, 4).Value = 44714
$ws.Cells.Item(275, 5).Value = 13
$ws.Cells.Item(275, 6).Value = 100112021
$ws.Cells.Item(275, 7).Value = "Ají"
$ws.Cells.Item(275, 8).Value = "Americana (o)"
$ws.Cells.Item(275, 9).Value = "Primera"
$ws.Cells.Item(275, 10).Value = 25
$ws.Cells.Item(275, 11).Value = 33000
$ws.Cells.Item(275, 12).Value = 34000
$ws.Cells.Item(275, 13).Value = 33480
$ws.Cells.Item(275, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(275, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(275, 16).Value = 1339
$ws.Cells.Item(275, 17).Value = 25
$ws.Cells.Item(275, 18).Value = "Hortaliza"

# --- New row 276 ---
$ws.Cells.Item(276, 1).Value = 9
$ws.Cells.Item(276, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(276, 3).Value = "Metropolitana"
$ws.Cells.Item(276, 4).Value = 44714
$ws.Cells.Item(276, 5).Value = 13
$ws.Cells.Item(276, 6).Value = 100112021
$ws.Cells.Item(276, 7).Value = "Ají"
$ws.Cells.Item(276, 8).Value = "Inferno"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 79
$ws.Cells.Item(276, 11).Value = 24000
$ws.Cells.Item(276, 12).Value = 25000
$ws.Cells.Item(276, 13).Value = 24494
$ws.Cells.Item(276, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(276, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(276, 16).Value = 2041
$ws.Cells.Item(276, 17).Value = 12
$ws.Cells.Item(276, 18).Value = "Hortaliza"
